$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Price": the whole database-generation for fuel prices was
# regrouped into one script, so the sheet is rebuilt with a new
# "Inland" column and refreshed Min/Avg/Max figures for Oil/Gas/Hrd
# plus new Moderate/Scarce figures for BIO/PEA.
# -----------------------------------------------------------------
$wsPrice = $wb.Worksheets("Price")
$wsPrice.Cells.Clear()

# Header row
$wsPrice.Range("B1").Value = "Domestic"
$wsPrice.Range("C1").Value = "Pipeline"
$wsPrice.Range("D1").Value = "Imported"
$wsPrice.Range("E1").Value = "Inland"
$wsPrice.Range("F1").Value = "Moderate"
$wsPrice.Range("G1").Value = "Scarce"

# Row labels
$wsPrice.Range("A2").Value = "Oil - Min"
$wsPrice.Range("A3").Value = "Oil - Avg"
$wsPrice.Range("A4").Value = "Oil - Max"
$wsPrice.Range("A5").Value = "Gas - Min"
$wsPrice.Range("A6").Value = "Gas - Avg"
$wsPrice.Range("A7").Value = "Gas - Max"
$wsPrice.Range("A8").Value = "Hrd - Min"
$wsPrice.Range("A9").Value = "Hrd - Avg"
$wsPrice.Range("A10").Value = "Hrd - Max"
$wsPrice.Range("A11").Value = "BIO"
$wsPrice.Range("A12").Value = "PEA"

# Pre-format the whole data block with the 2-decimal number style
# so that blank cells also carry the style like the reference file.
$wsPrice.Range("B2:G12").NumberFormat = "0.00"

# Oil - Avg (base) values
$wsPrice.Range("B3").Value = 29.79
$wsPrice.Range("C3").Value = 38.62
$wsPrice.Range("D3").Value = 53.32

# Gas - Avg (base) values
$wsPrice.Range("B6").Value = 19.55
$wsPrice.Range("C6").Value = 25.24
$wsPrice.Range("D6").Value = 33.52
$wsPrice.Range("E6").Value = 44.01

# Oil - Min = Oil - Avg * 0.9
$wsPrice.Range("B2").Formula = "=B3*0.9"
$wsPrice.Range("C2:D2").Formula = "=C3*0.9"

# Oil - Max = Oil - Avg * 1.1
$wsPrice.Range("B4").Formula = "=B3*1.1"
$wsPrice.Range("C4:D4").Formula = "=C3*1.1"

# Gas - Min = Gas - Avg * 0.9
$wsPrice.Range("B5").Formula = "=B6*0.9"
$wsPrice.Range("C5:E5").Formula = "=C6*0.9"

# Gas - Max = Gas - Avg * 1.1
$wsPrice.Range("B7").Formula = "=B6*1.1"
$wsPrice.Range("C7:E7").Formula = "=C6*1.1"

# Hrd Min/Avg/Max (Imported column only)
$wsPrice.Range("D8").Value = 75
$wsPrice.Range("D9").Value = 80
$wsPrice.Range("D10").Value = 90

# BIO / PEA (Moderate & Scarce columns only)
$wsPrice.Range("F11").Value = 10.08
$wsPrice.Range("G11").Value = 30.24
$wsPrice.Range("F12").Value = 9.36
$wsPrice.Range("G12").Value = 28.08

# -----------------------------------------------------------------
# Sheet "Status": three countries flip which gas-pipeline bucket
# (column G) they fall into; column I recomputes automatically
# through its existing formula.
# -----------------------------------------------------------------
$wsStatus = $wb.Worksheets("Status")
$wsStatus.Range("G19").Value = 1
$wsStatus.Range("G21").Value = 0
$wsStatus.Range("G22").Value = 0

# -----------------------------------------------------------------
# Sheet "Distance": the per-fuel (OIL/GAS/HRD) distance columns are
# grouped into a single column, refreshed with new figures for a
# few countries.
# -----------------------------------------------------------------
$wsDistance = $wb.Worksheets("Distance")
$wsDistance.Range("C1:D26").Delete()
$wsDistance.Range("B1").Value = "HRD"
$wsDistance.Range("B3").Value = 1508
$wsDistance.Range("B4").Value = 360
$wsDistance.Range("B19").Value = 1464
$wsDistance.Range("B23").Value = 1570

# -----------------------------------------------------------------
# Restore on-screen selections to match the saved view state.
# -----------------------------------------------------------------
$wsStatus.Activate()
$wsStatus.Range("G25").Select()

$wsDistance.Activate()
$wsDistance.Range("E9").Select()

$wsPrice.Activate()
$wsPrice.Range("E7").Select()
